$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend formatting (style) from the last existing data row (row 6) down to the new rows (7-9)
$ws.Range("A6:C6").Copy() | Out-Null
$ws.Range("A7:C9").PasteSpecial(-4122) | Out-Null

# Write header and all game rows (new slate of games replaces the old one)
$ws.Range("A1").Value = "NBA, Sunday 3rd Mar 2024"
$ws.Range("B1").Value = "Ballgorithm"
$ws.Range("C1").Value = "ESPN"

$ws.Range("A2").Value = "Philadelphia 76ers (34-25) vs Dallas Mavericks (34-26)"
$ws.Range("B2").Value = "Philadelphia 76ers (59.38%)"
$ws.Range("C2").Value = "Dallas Mavericks (72.4%)"

$ws.Range("A3").Value = "Golden State Warriors (32-27) vs Boston Celtics (47-12)"
$ws.Range("B3").Value = "Boston Celtics (90.32%)"
$ws.Range("C3").Value = "Boston Celtics (82.9%)"

$ws.Range("A4").Value = "Los Angeles Clippers (38-20) vs Minnesota Timberwolves (42-18)"
$ws.Range("B4").Value = "Minnesota Timberwolves (75.86%)"
$ws.Range("C4").Value = "Minnesota Timberwolves (63.6%)"

$ws.Range("A5").Value = "Detroit Pistons (9-50) vs Orlando Magic (34-26)"
$ws.Range("B5").Value = "Orlando Magic (71.43%)"
$ws.Range("C5").Value = "Orlando Magic (78.5%)"

$ws.Range("A6").Value = "Charlotte Hornets (15-45) vs Toronto Raptors (22-38)"
$ws.Range("B6").Value = "Charlotte Hornets (60.0%)"
$ws.Range("C6").Value = "Toronto Raptors (80.1%)"

$ws.Range("A7").Value = "New York Knicks (35-25) vs Cleveland Cavaliers (39-20)"
$ws.Range("B7").Value = "Cleveland Cavaliers (66.67%)"
$ws.Range("C7").Value = "Cleveland Cavaliers (73.9%)"

$ws.Range("A8").Value = "Indiana Pacers (34-28) vs San Antonio Spurs (12-48)"
$ws.Range("B8").Value = "Indiana Pacers (62.50%)"
$ws.Range("C8").Value = "Indiana Pacers (66.6%)"

$ws.Range("A9").Value = "Oklahoma City Thunder (41-18) vs Phoenix Suns (35-24)"
$ws.Range("B9").Value = "Oklahoma City Thunder (80.00%)"
$ws.Range("C9").Value = "Oklahoma City Thunder (64.1%)"

# Update column A width to fit the new (longer) text, matching the recorded bestFit width
$ws.Columns("A").ColumnWidth = 51.25

$ws.Range("A1").Select() | Out-Null
